# The edit inserts one new weekly price record as a new row 678 in the
# "Albahaca" (Hortaliza) sheet, pushing the former rows 678-769 down to
# 679-770 (dimension grows from A1:R769 to A1:R770).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 678; this shifts rows 678:769 down
# to 679:770 and Excel will recompute the sheet dimension automatically.
$ws.Rows("678:678").Insert()

# Populate the newly inserted row with the new record's data.
$ws.Cells.Item(678, 1).Value  = 6
$ws.Cells.Item(678, 2).Value  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(678, 3).Value  = "Metropolitana"
$ws.Cells.Item(678, 4).Value  = 45124
$ws.Cells.Item(678, 5).Value  = 13
$ws.Cells.Item(678, 6).Value  = 100112052
$ws.Cells.Item(678, 7).Value  = "Albahaca"
$ws.Cells.Item(678, 8).Value  = "Sin especificar"
$ws.Cells.Item(678, 9).Value  = "Primera"
$ws.Cells.Item(678, 10).Value = 45
$ws.Cells.Item(678, 11).Value = 4000
$ws.Cells.Item(678, 12).Value = 4500
$ws.Cells.Item(678, 13).Value = 4222
$ws.Cells.Item(678, 14).Value = "`$/paquete"
$ws.Cells.Item(678, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(678, 16).Value = 4222
$ws.Cells.Item(678, 17).Value = 1
$ws.Cells.Item(678, 18).Value = "Hortaliza"
